$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert 4 new blank rows right after row 38 (new rows become 39-42).
#    Row 38 itself (and everything above) keeps its position/content.
# ---------------------------------------------------------------------
$ws.Range("A39:A42").EntireRow.Insert(0) | Out-Null

# ---------------------------------------------------------------------
# 2. Update row 38: renumber, change description text, keep "DONE".
#    This also normalizes its formatting back to the plain table style
#    (no extra fill flag) to match the rest of the table rows.
# ---------------------------------------------------------------------
$ws.Cells.Item(38, 3).Value2 = 18.1
$ws.Cells.Item(38, 4).Value2 = "Employee and WorkingTimes filtering added"
$ws.Cells.Item(38, 5).Value2 = "DONE"
$ws.Range("C38:E38").Interior.Pattern = -4142   # xlNone - drop applyFill flag

# ---------------------------------------------------------------------
# 3. Fill in the 4 new rows (39-42) with their data, and give them the
#    same visual formatting as the rest of the table:
#    column C = bold numbering, columns C:E = thin box border.
# ---------------------------------------------------------------------
$rows = @(
    @{ Row = 39; Num = 18.2; Text = "Dinamičko mijenjanje linkova u header-u" },
    @{ Row = 40; Num = 19;   Text = "Forme za kreiranje company, employee i department" },
    @{ Row = 41; Num = 20.1; Text = "Create i Update za company" },
    @{ Row = 42; Num = 20.2; Text = "Create i Update za department" }
)

foreach ($item in $rows) {
    $r = $item.Row

    $cC = $ws.Cells.Item($r, 3)
    $cC.Value2 = $item.Num
    $cC.Font.Bold = $true
    $cC.Borders.LineStyle = 1

    $cD = $ws.Cells.Item($r, 4)
    $cD.Value2 = $item.Text
    $cD.Borders.LineStyle = 1

    $cE = $ws.Cells.Item($r, 5)
    $cE.Value2 = "DONE"
    $cE.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------
# 4. Match the view/selection state recorded in the saved workbook.
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("C41:E42").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
